# Trade #15 closed at 2026-02-17 07:59:04 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.92
$summary.Range("B4").Value = -0.08
$summary.Range("B5").Value = -0.11
$summary.Range("B6").Value = 15
$summary.Range("B8").Value = 8
$summary.Range("B9").Value = 33.33

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.92
$status.Range("D4").Value = 15
$status.Range("E4").Value = -0.08
$status.Range("F4").Value = -0.08
$status.Range("G4").Value = 33.33

# --- New trade row (#15) appended to "All Trades" and "MarketMaking" sheets ---
$newRow = @(15, "2026-02-17", "07:58:58", "MarketMaking", "UP", 0.08, 0.03, "CLOSED", -62.5, -0.05, 99.92, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    # Force the Date column (B) to remain plain text instead of being
    # auto-converted to a date serial number by Excel's type inference.
    $ws.Range("B16").NumberFormat = "@"
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(16, $i + 1).Value = $newRow[$i]
    }
}
